# generacion de iva total y neto
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New "Total" header column
$ws.Range("E1").Value = "Total"

# Net / IVA / Total rows
$ws.Range("E7").Formula = "=E11/1.19"

$ws.Range("E8").ClearContents()

$ws.Range("D10").Value = "IVA"
$ws.Range("E10").Formula = "=(E7*0.19)"

$ws.Range("D11").Value = "TOTAL"
$ws.Range("E11").Formula = "=SUM(E1:E6)"

# Workbook base/default font -> Century Gothic
$normal = $wb.Styles.Item(1)
$normal.Font.Name = "Century Gothic"

# Header row formatting: bold, size 14, Century Gothic
$ws.Range("A1:E1").Font.Bold = $true
$ws.Range("A1:E1").Font.Size = 14
$ws.Range("A1:E1").Font.Name = "Century Gothic"

$ws.Range("A1:E1").Select() | Out-Null
